$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (49) down into the new rows (50-56)
# so the new rows pick up the same column-A / column-E cell styles as the rest of the table.
$ws.Range("A49:V49").Copy()
$ws.Range("A50:V56").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 50 (Indice 49)
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "iran"
$ws.Cells.Item(50, 3).Value = "persian-gulf-pro-league"
$ws.Cells.Item(50, 4).Value = "2023-2024"
$ws.Cells.Item(50, 5).Value = 45226.6875
$ws.Cells.Item(50, 6).Value = "Esteghlal F.C."
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(50, 8).Value = "Aluminium Arak"
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 1.49
$ws.Cells.Item(50, 11).Value = "26/10/2023 07:42"
$ws.Cells.Item(50, 12).Value = 1.54
$ws.Cells.Item(50, 13).Value = "27/10/2023 16:24"
$ws.Cells.Item(50, 14).Value = 3.49
$ws.Cells.Item(50, 15).Value = "26/10/2023 07:42"
$ws.Cells.Item(50, 16).Value = 3.42
$ws.Cells.Item(50, 17).Value = "27/10/2023 16:25"
$ws.Cells.Item(50, 18).Value = 6.44
$ws.Cells.Item(50, 19).Value = "26/10/2023 07:42"
$ws.Cells.Item(50, 20).Value = 7.84
$ws.Cells.Item(50, 21).Value = "27/10/2023 16:25"
$ws.Cells.Item(50, 22).Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/esteghlal-teh-aluminium-arak/lzS4r9PG/"

# Row 51 (Indice 50)
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "iran"
$ws.Cells.Item(51, 3).Value = "persian-gulf-pro-league"
$ws.Cells.Item(51, 4).Value = "2023-2024"
$ws.Cells.Item(51, 5).Value = 45226.6875
$ws.Cells.Item(51, 6).Value = "Zob Ahan"
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = "Tractor"
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 3.1
$ws.Cells.Item(51, 11).Value = "26/10/2023 07:42"
$ws.Cells.Item(51, 12).Value = 3.13
$ws.Cells.Item(51, 13).Value = "27/10/2023 16:17"
$ws.Cells.Item(51, 14).Value = 2.74
$ws.Cells.Item(51, 15).Value = "26/10/2023 07:42"
$ws.Cells.Item(51, 16).Value = 2.69
$ws.Cells.Item(51, 17).Value = "27/10/2023 16:18"
$ws.Cells.Item(51, 18).Value = 2.35
$ws.Cells.Item(51, 19).Value = "26/10/2023 07:42"
$ws.Cells.Item(51, 20).Value = 2.63
$ws.Cells.Item(51, 21).Value = "27/10/2023 16:18"
$ws.Cells.Item(51, 22).Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/zob-ahan-tractor/I9W8sTvN/"

# Row 52 (Indice 51)
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = "iran"
$ws.Cells.Item(52, 3).Value = "persian-gulf-pro-league"
$ws.Cells.Item(52, 4).Value = "2023-2024"
$ws.Cells.Item(52, 5).Value = 45226.72916666666
$ws.Cells.Item(52, 6).Value = "Foolad"
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = "Shams Azar Qazvin"
$ws.Cells.Item(52, 9).Value = 1
$ws.Cells.Item(52, 10).Value = 1.74
$ws.Cells.Item(52, 11).Value = "26/10/2023 07:42"
$ws.Cells.Item(52, 12).Value = 1.95
$ws.Cells.Item(52, 13).Value = "27/10/2023 17:09"
$ws.Cells.Item(52, 14).Value = 2.99
$ws.Cells.Item(52, 15).Value = "26/10/2023 07:42"
$ws.Cells.Item(52, 16).Value = 2.79
$ws.Cells.Item(52, 17).Value = "27/10/2023 17:09"
$ws.Cells.Item(52, 18).Value = 4.88
$ws.Cells.Item(52, 19).Value = "26/10/2023 07:42"
$ws.Cells.Item(52, 20).Value = 5.01
$ws.Cells.Item(52, 21).Value = "27/10/2023 17:09"
$ws.Cells.Item(52, 22).Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/foolad-shams-azar-qazvin/d0VCtmfT/"

# Row 53 (Indice 52)
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = "iran"
$ws.Cells.Item(53, 3).Value = "persian-gulf-pro-league"
$ws.Cells.Item(53, 4).Value = "2023-2024"
$ws.Cells.Item(53, 5).Value = 45227.66666666666
$ws.Cells.Item(53, 6).Value = "Gol Gohar"
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = "Havadar SC"
$ws.Cells.Item(53, 9).Value = 3
$ws.Cells.Item(53, 10).Value = 1.98
$ws.Cells.Item(53, 11).Value = "27/10/2023 13:43"
$ws.Cells.Item(53, 12).Value = 1.92
$ws.Cells.Item(53, 13).Value = "28/10/2023 15:55"
$ws.Cells.Item(53, 14).Value = 2.85
$ws.Cells.Item(53, 15).Value = "27/10/2023 13:43"
$ws.Cells.Item(53, 16).Value = 2.68
$ws.Cells.Item(53, 17).Value = "28/10/2023 15:55"
$ws.Cells.Item(53, 18).Value = 3.99
$ws.Cells.Item(53, 19).Value = "27/10/2023 13:43"
$ws.Cells.Item(53, 20).Value = 5.64
$ws.Cells.Item(53, 21).Value = "28/10/2023 15:55"
$ws.Cells.Item(53, 22).Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/gol-gohar-havadar-sc/2s8ux72p/"

# Row 54 (Indice 53)
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = "iran"
$ws.Cells.Item(54, 3).Value = "persian-gulf-pro-league"
$ws.Cells.Item(54, 4).Value = "2023-2024"
$ws.Cells.Item(54, 5).Value = 45227.6875
$ws.Cells.Item(54, 6).Value = "Paykan"
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = "Mes Rafsanjan"
$ws.Cells.Item(54, 9).Value = 1
$ws.Cells.Item(54, 10).Value = 2.98
$ws.Cells.Item(54, 11).Value = "27/10/2023 13:43"
$ws.Cells.Item(54, 12).Value = 4.52
$ws.Cells.Item(54, 13).Value = "28/10/2023 16:10"
$ws.Cells.Item(54, 14).Value = 2.58
$ws.Cells.Item(54, 15).Value = "27/10/2023 13:43"
$ws.Cells.Item(54, 16).Value = 2.44
$ws.Cells.Item(54, 17).Value = "28/10/2023 16:10"
$ws.Cells.Item(54, 18).Value = 2.57
$ws.Cells.Item(54, 19).Value = "27/10/2023 13:43"
$ws.Cells.Item(54, 20).Value = 2.27
$ws.Cells.Item(54, 21).Value = "28/10/2023 16:10"
$ws.Cells.Item(54, 22).Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/paykan-mes-rafsanjan/GpCqyRHj/"

# Row 55 (Indice 54)
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = "iran"
$ws.Cells.Item(55, 3).Value = "persian-gulf-pro-league"
$ws.Cells.Item(55, 4).Value = "2023-2024"
$ws.Cells.Item(55, 5).Value = 45228.52083333334
$ws.Cells.Item(55, 6).Value = "Nassaji Mazandaran"
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = "Sepahan"
$ws.Cells.Item(55, 9).Value = 1
$ws.Cells.Item(55, 10).Value = 7.17
$ws.Cells.Item(55, 11).Value = "29/10/2023 10:42"
$ws.Cells.Item(55, 12).Value = 7.27
$ws.Cells.Item(55, 13).Value = "29/10/2023 12:06"
$ws.Cells.Item(55, 14).Value = 4.03
$ws.Cells.Item(55, 15).Value = "29/10/2023 10:42"
$ws.Cells.Item(55, 16).Value = 3.91
$ws.Cells.Item(55, 17).Value = "29/10/2023 12:06"
$ws.Cells.Item(55, 18).Value = 1.46
$ws.Cells.Item(55, 19).Value = "29/10/2023 10:42"
$ws.Cells.Item(55, 20).Value = 1.48
$ws.Cells.Item(55, 21).Value = "29/10/2023 12:06"
$ws.Cells.Item(55, 22).Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mazandaran-sepahan/8vaVbU9c/"

# Row 56 (Indice 55)
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = "iran"
$ws.Cells.Item(56, 3).Value = "persian-gulf-pro-league"
$ws.Cells.Item(56, 4).Value = "2023-2024"
$ws.Cells.Item(56, 5).Value = 45228.5625
$ws.Cells.Item(56, 6).Value = "Malavan"
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = "Persepolis"
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 6.19
$ws.Cells.Item(56, 11).Value = "28/10/2023 18:13"
$ws.Cells.Item(56, 12).Value = 7.22
$ws.Cells.Item(56, 13).Value = "29/10/2023 12:04"
$ws.Cells.Item(56, 14).Value = 3.63
$ws.Cells.Item(56, 15).Value = "28/10/2023 18:13"
$ws.Cells.Item(56, 16).Value = 3.6
$ws.Cells.Item(56, 17).Value = "29/10/2023 12:04"
$ws.Cells.Item(56, 18).Value = 1.53
$ws.Cells.Item(56, 19).Value = "28/10/2023 18:13"
$ws.Cells.Item(56, 20).Value = 1.53
$ws.Cells.Item(56, 21).Value = "29/10/2023 09:40"
$ws.Cells.Item(56, 22).Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/malavan-persepolis/ljGSclP3/"

